# Add "2022-Q4" sheet (new quarterly holdings snapshot) right after "总计",
# matching the commit "feat: add 2022-Q4 data".
#
# Effects:
#  1. New worksheet "2022-Q4" inserted immediately after "总计" (shifts all
#     other quarter sheets one position to the right; their sheetIds/parts
#     are renumbered automatically by the host on save).
#  2. "总计" (summary) sheet gets a new row 2 with the 2022-Q4 totals; the
#     existing rows shift down one position and keep their fixed row-index
#     (column A) numbering (0,1,2,3,4,5) rather than carrying their old
#     index with them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing text storage (keeps
# leading zeros / decimal-looking strings as literal text instead of
# Excel auto-coercing them to numbers), and without leaving a stray
# cell style behind.
# ---------------------------------------------------------------------
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)

$newSheet = $wb.Worksheets.Add($null, $summarySheet)
$newSheet.Name = "2022-Q4"

# Template cells (on the existing "2022-Q3" sheet) used purely to copy
# the already-registered cell style ("bold + thin border + center/top",
# style index 2 in the original workbook) onto the new sheet's header
# row and index column, so we reuse the existing style instead of
# minting new ones.
$headerStyleTemplate = $q3Sheet.Range("B1")
$indexStyleTemplate = $q3Sheet.Range("A2")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}
$headerStyleTemplate.Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Row data: idx, code, name, scale(基金规模), position(股票总仓位),
# ratio(仓位占比), marketValue(持有市值), rank(仓位排名)
$q4Data = @(
        @(0, '003835', '鹏华沪深港新兴成长灵活配置混合', '45.58', '94.28', '9.04', '4.1204', 1),
        @(1, '001018', '易方达新经济灵活配置混合', '74.11', '93.12', '4.94', '3.6610', 3),
        @(2, '110013', '易方达科翔混合', '67.77', '92.12', '4.84', '3.2801', 3),
        @(3, '009341', '易方达均衡成长股票', '59.89', '88.68', '4.24', '2.5393', 4),
        @(4, '501203', '易方达创新未来混合（LOF）', '52.13', '89.00', '4.43', '2.3094', 4),
        @(5, '000729', '建信中小盘先锋股票A', '31.93', '91.13', '5.29', '1.6891', 2),
        @(6, '016067', '鹏华新能源汽车混合A', '19.04', '95.25', '8.84', '1.6831', 1),
        @(7, '110001', '易方达平稳增长混合', '35.60', '62.38', '3.75', '1.3350', 4),
        @(8, '001076', '易方达改革红利混合', '25.93', '93.28', '4.92', '1.2758', 3),
        @(9, '530005', '建信优化配置混合A', '24.05', '88.67', '5.06', '1.2169', 2),
        @(10, '000756', '建信潜力新蓝筹股票A', '15.30', '90.58', '5.02', '0.7681', 2),
        @(11, '013919', '建信中小盘先锋股票C', '10.14', '91.13', '5.29', '0.5364', 2),
        @(12, '014967', '建信潜力新蓝筹股票C', '10.68', '90.58', '5.02', '0.5361', 2),
        @(13, '002418', '汇添富优选回报灵活配置混合C', '9.25', '94.49', '4.81', '0.4449', 9),
        @(14, '001166', '建信环保产业股票', '7.08', '86.02', '4.92', '0.3483', 5),
        @(15, '016068', '鹏华新能源汽车混合C', '3.65', '95.25', '8.84', '0.3227', 1),
        @(16, '470021', '汇添富优选回报灵活配置混合A', '5.37', '94.49', '4.81', '0.2583', 9),
        @(17, '011460', '鹏华创新成长混合A', '10.08', '78.91', '2.43', '0.2449', 8),
        @(18, '008786', '长城健康生活灵活配置混合', '5.69', '83.10', '3.05', '0.1735', 5),
        @(19, '004671', '中融核心成长灵活配置混合', '1.14', '65.86', '4.00', '0.0456', 7),
        @(20, '010009', '中融成长优选混合C', '1.05', '60.81', '3.81', '0.0400', 6),
        @(21, '159804', '国寿安保国证创业板中盘精选88ETF', '1.15', '99.00', '2.42', '0.0278', 1),
        @(22, '010008', '中融成长优选混合A', '0.57', '60.81', '3.81', '0.0217', 6),
        @(23, '004536', '嘉实中小企业量化活力灵活配置混合', '0.22', '93.77', '5.61', '0.0123', 2),
        @(24, '011461', '鹏华创新成长混合C', '0.35', '78.91', '2.43', '0.0085', 8),
        @(25, '015436', '建信优化配置混合C', '0.11', '88.67', '5.06', '0.0056', 2),
        @(26, '006538', '东海核心价值精选混合', '0.12', '85.21', '2.73', '0.0033', 6),
        @(27, '005281', '中科沃土转型升级灵活配置混合', '0.10', '60.18', '2.75', '0.0028', 9)
)

foreach ($row in $q4Data) {
    $r = [int]$row[0] + 2

    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $indexStyleTemplate.Copy()
    $newSheet.Cells.Item($r, 1).PasteSpecial(-4122)
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    Set-TextValue $newSheet.Cells.Item($r, 2) $row[1]
    Set-TextValue $newSheet.Cells.Item($r, 3) $row[2]
    Set-TextValue $newSheet.Cells.Item($r, 4) $row[3]
    Set-TextValue $newSheet.Cells.Item($r, 5) $row[4]
    Set-TextValue $newSheet.Cells.Item($r, 6) $row[5]
    Set-TextValue $newSheet.Cells.Item($r, 7) $row[6]

    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2) Insert a new row into "总计" for the 2022-Q4 totals, then fix up
#    the fixed row-index column (A) for every data row (0,1,2,3,4,5).
# ---------------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()
$summarySheet.Rows.Item(2).ClearFormats()

$summarySheet.Range("B2").Value = "2022-Q4"
$summarySheet.Range("C2").Value = 28
$summarySheet.Range("D2").Value = 26.91

for ($r = 2; $r -le 7; $r++) {
    $cell = $summarySheet.Cells.Item($r, 1)
    $indexStyleTemplate.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $r - 2
}

# ---------------------------------------------------------------------
# 3) Restore "总计" as the active sheet/tab (matches unchanged
#    workbook-level bookView activeTab="0").
# ---------------------------------------------------------------------
$summarySheet.Activate()
$summarySheet.Range("A1").Select()
